$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 4).Value = '60.731.58'
$ws.Cells.Item(2, 5).Value = '  +2.13%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.689.89'
$ws.Cells.Item(3, 5).Value = '  +1.84%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.37%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '521.35'
$ws.Cells.Item(5, 5).Value = '  +0.80%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '149.13'
$ws.Cells.Item(6, 5).Value = '  +1.26%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.24%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.580'
$ws.Cells.Item(8, 5).Value = '  +1.32%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.708.45'
$ws.Cells.Item(9, 5).Value = '  +1.41%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '6.43'
$ws.Cells.Item(10, 5).Value = '  -0.50%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.53%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.86%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +1.14%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '3.154.55'
$ws.Cells.Item(14, 5).Value = '  +0.70%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '60.698.03'
$ws.Cells.Item(15, 5).Value = '  +2.09%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '21.51'
$ws.Cells.Item(16, 5).Value = '  +1.13%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '2.833.46'
$ws.Cells.Item(17, 5).Value = '  +6.21%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).Value = '0.0000140'
$ws.Cells.Item(18, 5).Value = '  +0.95%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '354.68'
$ws.Cells.Item(19, 5).Value = '  +2.48%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '4.59'
$ws.Cells.Item(20, 5).Value = '  -0.55%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '10.58'
$ws.Cells.Item(21, 5).Value = '  +0.46%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '6.33'
$ws.Cells.Item(22, 5).Value = '  +1.92%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.02%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '62.92'
$ws.Cells.Item(24, 5).Value = '  +2.34%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '0.424'
$ws.Cells.Item(25, 5).Value = '  +0.07%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +3.09%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -0.45%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '0.0₃0832'
$ws.Cells.Item(28, 5).Value = '  +1.17%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '7.34'
$ws.Cells.Item(29, 5).Value = '  +1.88%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '6.80'
$ws.Cells.Item(30, 5).Value = '  +4.33%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -0.21%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '19.19'
$ws.Cells.Item(32, 5).Value = '  +0.40%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '1.60'
$ws.Cells.Item(33, 5).Value = '  +0.67%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '150.48'
$ws.Cells.Item(34, 5).Value = '  +0.33%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +2.52%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '0.969'
$ws.Cells.Item(36, 5).Value = '  -7.50%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +3.66%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '1.57'
$ws.Cells.Item(38, 5).Value = '  +9.49%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.887'
$ws.Cells.Item(39, 5).Value = '  +1.71%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '36.90'
$ws.Cells.Item(40, 5).Value = '  +0.61%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '3.75'
$ws.Cells.Item(41, 5).Value = '  +0.34%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '287.42'
$ws.Cells.Item(42, 5).Value = '  +0.29%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'EnergySwap'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(43, 4).Value = '20.16'
$ws.Cells.Item(43, 5).Value = '  +1.78%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Stellar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(44, 4).Value = '0.0993'
$ws.Cells.Item(44, 5).Value = '  -0.21%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '0.612'
$ws.Cells.Item(45, 5).Value = '  -1.18%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.01%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '2.115.75'
$ws.Cells.Item(47, 5).Value = '  +5.97%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +5.43%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '0.0541'
$ws.Cells.Item(49, 5).Value = '  -0.76%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '0.0235'
$ws.Cells.Item(50, 5).Value = '  +0.56%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '19.29'
$ws.Cells.Item(51, 5).Value = '  +3.78%  '
